# crosscheck.xlsx: "added all scripts for all file types"
#
# - H column (year header 2010 -> 2012) no longer carries the per-row "X"
#   script-complete marks; those are cleared out (H3:H29).
# - J column ("OK GENERADOR") gets the remaining "X" marks filled in
#   (J16, J25:J29 previously blank) and every mark in J8:J23 / J25:J29 is
#   highlighted with a green fill to show the generator script now covers
#   every file type.
# - A new note is added below the existing one explaining 2012 already has
#   the full 24 departments.
# - The view is scrolled down a bit and the last selected cell updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header: year column H goes from 2010 to 2012 --------------------
$ws.Range("H2").Value = 2012

# ---- Clear the old "X" marks out of column H (rows 3-29) --------------
$ws.Range("H3:H29").ClearContents()

# ---- Fill in the missing "X" marks in column J and mark them ----------
# J16 was blank before; J25:J29 were blank before. All of J8:J23 and
# J25:J29 get the new green highlight fill applied (RGB 146,208,80 /
# hex 92D050), matching the newly-added xf that turns on fillId=2.
$green = 5296274  # OLE BGR packed value for RGB(146,208,80) = #92D050

$ws.Range("J16").Value = "X"
$ws.Range("J25").Value = "X"
$ws.Range("J26").Value = "X"
$ws.Range("J27").Value = "X"
$ws.Range("J28").Value = "X"
$ws.Range("J29").Value = "X"

$ws.Range("J8:J23").Interior.Color = $green
$ws.Range("J25:J29").Interior.Color = $green

# ---- New note under the existing 2007/CIIU note ------------------------
$ws.Range("E34").Value = "Nota: 2012 tiene ya los 24 deptos"

# ---- Update the view: scroll down/right a bit and move the selection --
$ws.Range("E35").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 3
